$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("contacts")
$ws.Activate()

$ws.Range("A5").Value = "Dr."
$ws.Range("B5").Value = "Peter"
$ws.Range("C5").Value = "Cris"
$ws.Range("D5").Value = "Flipkart"

$ws.Range("D5").Select()
